$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for B1
$ws.Range("B1").Value = "first_release_value"

# Apply the date style (from A2) to the new date cells A3:A22 first
$ws.Range("A2").Copy()
$ws.Range("A3:A22").PasteSpecial(-4122)

# Write date + value data
$ws.Cells.Item(2, 1).Value = 38717
$ws.Cells.Item(2, 2).Value = $null

$ws.Cells.Item(3, 1).Value = 39082
$ws.Cells.Item(3, 2).Value = -0.1287148287979267

$ws.Cells.Item(4, 1).Value = 39447
$ws.Cells.Item(4, 2).Value = 0.08844991283951664

$ws.Cells.Item(5, 1).Value = 39813
$ws.Cells.Item(5, 2).Value = 0.9225722794137248

$ws.Cells.Item(6, 1).Value = 40178
$ws.Cells.Item(6, 2).Value = -0.7475385776494314

$ws.Cells.Item(7, 1).Value = 40543
$ws.Cells.Item(7, 2).Value = 1.094490700431927

$ws.Cells.Item(8, 1).Value = 40908
$ws.Cells.Item(8, 2).Value = -0.8742167833903691

$ws.Cells.Item(9, 1).Value = 41274
$ws.Cells.Item(9, 2).Value = -0.5751029748885195

$ws.Cells.Item(10, 1).Value = 41639
$ws.Cells.Item(10, 2).Value = -0.1263097576649996

$ws.Cells.Item(11, 1).Value = 42004
$ws.Cells.Item(11, 2).Value = 0.4976690624053814

$ws.Cells.Item(12, 1).Value = 42369
$ws.Cells.Item(12, 2).Value = -0.3267144271395628

$ws.Cells.Item(13, 1).Value = 42735
$ws.Cells.Item(13, 2).Value = 0.149524011641855

$ws.Cells.Item(14, 1).Value = 43100
$ws.Cells.Item(14, 2).Value = 0.1983963998054783

$ws.Cells.Item(15, 1).Value = 43465
$ws.Cells.Item(15, 2).Value = 0.6497679376401333

$ws.Cells.Item(16, 1).Value = 43830
$ws.Cells.Item(16, 2).Value = -0.5803176690338252

$ws.Cells.Item(17, 1).Value = 44196
$ws.Cells.Item(17, 2).Value = 0.2929419500579789

$ws.Cells.Item(18, 1).Value = 44561
$ws.Cells.Item(18, 2).Value = -2.411738983600742

$ws.Cells.Item(19, 1).Value = 44926
$ws.Cells.Item(19, 2).Value = -1.022826494952023

$ws.Cells.Item(20, 1).Value = 45291
$ws.Cells.Item(20, 2).Value = -0.5858537819409149

$ws.Cells.Item(21, 1).Value = 45657
$ws.Cells.Item(21, 2).Value = 0.01773204329378331

$ws.Cells.Item(22, 1).Value = 46022
$ws.Cells.Item(22, 2).Value = $null
